$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '43.318.68'
$ws.Range("E2").Value = '  -4.68%  '
$ws.Range("D3").Value = '2.244.13'
$ws.Range("E3").Value = '  -5.62%  '
$ws.Range("E4").Value = '  -0.17%  '
$ws.Range("D5").Value = "'321.05"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +1.50%  '
$ws.Range("D6").Value = "'101.72"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  -6.57%  '
$ws.Range("D7").Value = "'0.588"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  -7.89%  '
$ws.Range("E8").Value = '  -0.20%  '
$ws.Range("D9").Value = "'0.567"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  -8.01%  '
$ws.Range("D10").Value = "'37.25"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  -9.13%  '
$ws.Range("D11").Value = "'54.56"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  -2.48%  '
$ws.Range("D12").Value = "'0.0831"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  -9.78%  '
$ws.Range("D13").Value = "'7.75"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  -9.46%  '
$ws.Range("E14").Value = '  -0.85%  '
$ws.Range("B15").Value = 'Polygon'
$ws.Range("C15").Value = 'https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic'
$ws.Range("D15").Value = "'0.871"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  -11.69%  '
$ws.Range("B16").Value = 'WrappedliquidstakedEther2.0'
$ws.Range("C16").Value = 'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth'
$ws.Range("D16").Value = '2.577.95'
$ws.Range("E16").Value = '  -5.85%  '
$ws.Range("D17").Value = "'14.48"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  -6.66%  '
$ws.Range("D18").Value = '2.236.70'
$ws.Range("E18").Value = '  -5.69%  '
$ws.Range("D19").Value = '43.212.09'
$ws.Range("E19").Value = '  -4.84%  '
$ws.Range("D20").Value = "'14.53"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  -10.35%  '
$ws.Range("D21").Value = '0.0₃0971'
$ws.Range("E21").Value = '  -9.08%  '
$ws.Range("D22").Value = "'6.55"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  -10.91%  '
$ws.Range("D23").Value = "'65.71"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  -10.41%  '
$ws.Range("D24").Value = "'3.22"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  -11.64%  '
$ws.Range("D25").Value = "'238.59"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  -8.73%  '
$ws.Range("D26").Value = "'2.17"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  -8.13%  '
$ws.Range("E27").Value = '  -0.20%  '
$ws.Range("E28").Value = '  +2.13%  '
$ws.Range("D29").Value = "'10.09"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  -9.96%  '
$ws.Range("D30").Value = "'2.24"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  -2.44%  '
$ws.Range("D31").Value = "'6.42"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  -16.37%  '
$ws.Range("D32").Value = "'35.86"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  -3.83%  '
$ws.Range("D33").Value = "'20.54"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  -8.39%  '
$ws.Range("D34").Value = "'0.0879"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  -9.02%  '
$ws.Range("D35").Value = "'154.03"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  -7.73%  '
$ws.Range("D36").Value = "'2.74"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  -4.77%  '
$ws.Range("D37").Value = "'3.18"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  +8.44%  '
$ws.Range("D38").Value = "'1.98"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  +3.46%  '
$ws.Range("E39").Value = '  -7.74%  '
$ws.Range("D40").Value = "'4.47"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  -5.42%  '
$ws.Range("D41").Value = "'0.104"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  -11.38%  '
$ws.Range("D42").Value = "'3.71"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  -9.23%  '
$ws.Range("D43").Value = "'0.0326"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  -8.52%  '
$ws.Range("D44").Value = "'13.16"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  -0.08%  '
$ws.Range("E45").Value = '  -0.12%  '
$ws.Range("D46").Value = '1.782.70'
$ws.Range("E46").Value = '  -2.44%  '
$ws.Range("D47").Value = "'86.92"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  -11.42%  '
$ws.Range("D48").Value = "'0.207"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  -9.92%  '
$ws.Range("D49").Value = "'5.34"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  -11.55%  '
$ws.Range("D50").Value = "'76.41"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  -9.63%  '
$ws.Range("D51").Value = "'59.51"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  -15.82%  '
